# Auto-generated script applying market-data refresh per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 443.77777
$ws.Cells.Item(28, 9).Value = 374.25
$ws.Cells.Item(28, 11).Value = 374.25
$ws.Cells.Item(28, 13).Value = 110.75
$ws.Cells.Item(33, 8).Value = 304.66666
$ws.Cells.Item(33, 9).Value = 223.8421
$ws.Cells.Item(33, 11).Value = 223.8421
$ws.Cells.Item(33, 13).Value = 5.157900000000012
$ws.Cells.Item(39, 8).Value = 272.2
$ws.Cells.Item(39, 9).Value = 194.16667
$ws.Cells.Item(39, 10).Value = 389.25
$ws.Cells.Item(39, 11).Value = 582.50001
$ws.Cells.Item(39, 12).Value = 1167.75
$ws.Cells.Item(39, 13).Value = -286.50001
$ws.Cells.Item(39, 14).Value = -1759.75
$ws.Cells.Item(40, 8).Value = 1851.5
$ws.Cells.Item(40, 9).Value = 2002
$ws.Cells.Item(40, 10).Value = 1801.3334
$ws.Cells.Item(40, 11).Value = 2002
$ws.Cells.Item(40, 12).Value = 1801.3334
$ws.Cells.Item(40, 13).Value = -1827
$ws.Cells.Item(40, 14).Value = -2151.3334
$ws.Cells.Item(80, 8).Value = 1014.2857
$ws.Cells.Item(80, 9).Value = 1000
$ws.Cells.Item(80, 10).Value = 1016.6667
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 12).Value = 3050.0001
$ws.Cells.Item(80, 13).Value = -2002
$ws.Cells.Item(80, 14).Value = -5046.0001
$ws.Cells.Item(83, 8).Value = 1014.2857
$ws.Cells.Item(83, 9).Value = 1000
$ws.Cells.Item(83, 10).Value = 1016.6667
$ws.Cells.Item(83, 11).Value = 9000
$ws.Cells.Item(83, 12).Value = 9150.0003
$ws.Cells.Item(83, 13).Value = -4008
$ws.Cells.Item(83, 14).Value = -19134.0003
$ws.Cells.Item(113, 8).Value = 4333
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 4999.5
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 4999.5
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(113, 14).Value = -11507.5
$ws.Cells.Item(137, 8).Value = 2898.75
$ws.Cells.Item(137, 9).Value = 2898.75
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 8696.25
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -6146.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12353.4
$ws.Cells.Item(32, 10).Value = 24999.75
$ws.Cells.Item(32, 12).Value = 24999.75
$ws.Cells.Item(32, 14).Value = -25573.75
$ws.Cells.Item(61, 8).Value = 23494.5
$ws.Cells.Item(61, 10).Value = 4990
$ws.Cells.Item(61, 12).Value = 4990
$ws.Cells.Item(61, 14).Value = -5414
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 2627.7896
$ws.Cells.Item(132, 9).Value = 2627.7896
$ws.Cells.Item(132, 11).Value = 7883.3688
$ws.Cells.Item(132, 13).Value = -5353.3688
$ws.Cells.Item(136, 8).Value = 23494.5
$ws.Cells.Item(136, 10).Value = 4990
$ws.Cells.Item(136, 12).Value = 14970
$ws.Cells.Item(136, 14).Value = -20070
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2085.125
$ws.Cells.Item(134, 9).Value = 2085.125
$ws.Cells.Item(134, 11).Value = 6255.375
$ws.Cells.Item(134, 13).Value = -3720.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3060.8948
$ws.Cells.Item(31, 9).Value = 3073.5
$ws.Cells.Item(31, 11).Value = 3073.5
$ws.Cells.Item(31, 13).Value = -2778.5
$ws.Cells.Item(34, 8).Value = 3060.8948
$ws.Cells.Item(34, 9).Value = 3073.5
$ws.Cells.Item(34, 11).Value = 3073.5
$ws.Cells.Item(34, 13).Value = -2871.5
$ws.Cells.Item(58, 8).Value = 2756
$ws.Cells.Item(58, 9).Value = 2756
$ws.Cells.Item(58, 11).Value = 2756
$ws.Cells.Item(58, 13).Value = -2553
$ws.Cells.Item(122, 8).Value = 3999
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(136, 8).Value = 2756
$ws.Cells.Item(136, 9).Value = 2756
$ws.Cells.Item(136, 11).Value = 8268
$ws.Cells.Item(136, 13).Value = -5718
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 1457.1428
$ws.Cells.Item(55, 10).Value = 1500
$ws.Cells.Item(55, 12).Value = 4500
$ws.Cells.Item(55, 14).Value = -4854
$ws.Cells.Item(57, 8).Value = 3000
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 2003715.4
$ws.Cells.Item(3, 9).Value = 1337667.6
$ws.Cells.Item(3, 10).Value = 2503251.2
$ws.Cells.Item(3, 11).Value = 1337667.6
$ws.Cells.Item(3, 12).Value = 2503251.2
$ws.Cells.Item(3, 13).Value = -1337551.6
$ws.Cells.Item(3, 14).Value = -2503483.2
$ws.Cells.Item(54, 8).Value = 18500
$ws.Cells.Item(54, 10).Value = 18500
$ws.Cells.Item(54, 12).Value = 18500
$ws.Cells.Item(54, 14).Value = -19280
$ws.Cells.Item(80, 8).Value = 13990
$ws.Cells.Item(80, 9).Value = 3683.3333
$ws.Cells.Item(80, 10).Value = 29450
$ws.Cells.Item(80, 11).Value = 3683.3333
$ws.Cells.Item(80, 12).Value = 29450
$ws.Cells.Item(80, 13).Value = -2685.3333
$ws.Cells.Item(80, 14).Value = -31446
$ws.Cells.Item(83, 8).Value = 13990
$ws.Cells.Item(83, 9).Value = 3683.3333
$ws.Cells.Item(83, 10).Value = 29450
$ws.Cells.Item(83, 11).Value = 18416.6665
$ws.Cells.Item(83, 12).Value = 147250
$ws.Cells.Item(83, 13).Value = -13424.6665
$ws.Cells.Item(83, 14).Value = -157234
$ws.Cells.Item(113, 8).Value = 1725
$ws.Cells.Item(113, 9).Value = 1725
$ws.Cells.Item(113, 11).Value = 1725
$ws.Cells.Item(113, 13).Value = 445
$ws.Cells.Item(117, 8).Value = 65000
$ws.Cells.Item(117, 10).Value = 65000
$ws.Cells.Item(117, 12).Value = 65000
$ws.Cells.Item(117, 14).Value = -71884
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 2888.0833
$ws.Cells.Item(132, 10).Value = 3675
$ws.Cells.Item(132, 12).Value = 11025
$ws.Cells.Item(132, 14).Value = -16085
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1600
$ws.Cells.Item(7, 9).Value = 1600
$ws.Cells.Item(7, 11).Value = 1600
$ws.Cells.Item(7, 13).Value = -1488
$ws.Cells.Item(16, 8).Value = 7218.75
$ws.Cells.Item(16, 9).Value = 7218.75
$ws.Cells.Item(16, 11).Value = 7218.75
$ws.Cells.Item(16, 13).Value = -7048.75
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).ClearContents()
$ws.Cells.Item(125, 14).Value = 0
$ws.Cells.Item(126, 8).Value = 1600
$ws.Cells.Item(126, 9).Value = 1600
$ws.Cells.Item(126, 11).Value = 4800
$ws.Cells.Item(126, 13).Value = -2330
$ws.Cells.Item(132, 8).Value = 5236.8887
$ws.Cells.Item(132, 9).Value = 3826.2
$ws.Cells.Item(132, 11).Value = 11478.6
$ws.Cells.Item(132, 13).Value = -8948.599999999999
$ws.Cells.Item(136, 8).Value = 94469
$ws.Cells.Item(136, 9).Value = 35161.6
$ws.Cells.Item(136, 10).Value = 193314.67
$ws.Cells.Item(136, 11).Value = 105484.8
$ws.Cells.Item(136, 12).Value = 579944.01
$ws.Cells.Item(136, 13).Value = -102934.8
$ws.Cells.Item(136, 14).Value = -585044.01
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 1000
$ws.Cells.Item(86, 10).Value = 1000
$ws.Cells.Item(86, 12).Value = 1000
$ws.Cells.Item(86, 14).Value = -3246
$ws.Cells.Item(89, 8).Value = 1000
$ws.Cells.Item(89, 10).Value = 1000
$ws.Cells.Item(89, 12).Value = 5000
$ws.Cells.Item(89, 14).Value = -16232
$ws.Cells.Item(93, 8).Value = 40000
$ws.Cells.Item(93, 10).Value = 40000
$ws.Cells.Item(93, 12).Value = 40000
$ws.Cells.Item(93, 14).Value = -44992
$ws.Cells.Item(122, 8).Value = 10880.6
$ws.Cells.Item(122, 9).Value = 1100.75
$ws.Cells.Item(122, 11).Value = 3302.25
$ws.Cells.Item(122, 13).Value = -852.25
$ws.Cells.Item(132, 8).Value = 1463.1818
$ws.Cells.Item(132, 9).Value = 1237
$ws.Cells.Item(132, 10).Value = 2066.3333
$ws.Cells.Item(132, 11).Value = 3711
$ws.Cells.Item(132, 12).Value = 6198.999899999999
$ws.Cells.Item(132, 13).Value = -1181
$ws.Cells.Item(132, 14).Value = -11258.9999
$ws.Cells.Item(136, 8).Value = 4500
$ws.Cells.Item(136, 9).Value = 4500
$ws.Cells.Item(136, 11).Value = 13500
$ws.Cells.Item(136, 13).Value = -10950

Write-Host "Applied 201 cell updates across 8 sheets"
